$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 2.72
$ws.Range("I4").Value = 2.55
$ws.Range("W4").Value = 8.75
$ws.Range("X4").Value = 14
$ws.Range("AC4").Value = 9
$ws.Range("AD4").Value = 5.9
$ws.Range("AH4").Value = 8.5
$ws.Range("AJ4").Value = 9.25
$ws.Range("AN4").Value = 4.7
$ws.Range("AO4").Value = 14.5
$ws.Range("AP4").Value = 19.5
$ws.Range("AU4").Value = 6.3
